$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the SUM ranges in B24/C24 (and consequently A24/A25) to include
# rows 11-18 instead of rows 13-18
$ws.Range("B24").Formula = "=SUM(B11:B18)"
$ws.Range("C24").Formula = "=SUM(C11:C18)"

# Bump the row height for the whole used area (1:25) from 15 to 15.75
$ws.Rows("1:25").RowHeight = 15.75

# Update the selection to match the new analysis range (B11:C18), active cell B11
$ws.Range("B11:C18").Select()
